$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_30_2_24"
$ws.Range("B2").Value = 0.9998858021297656
$ws.Range("C2").Value = 0.9992155822822145
$ws.Range("D2").Value = 0.9993485757971506
$ws.Range("E2").Value = 0.9998294088637025
$ws.Range("F2").Value = 0.9996432421476592
$ws.Range("G2").Value = 0.0001065986991565944
$ws.Range("H2").Value = 0.0007322195075944053
$ws.Range("I2").Value = 0.0004848459151633974
$ws.Range("J2").Value = 0.0001802514907205349
$ws.Range("K2").Value = 0.0003325487029419661
$ws.Range("L2").Value = 0.000628115785851762
$ws.Range("M2").Value = 0.01032466460262
$ws.Range("N2").Value = 1.000109629955425
$ws.Range("O2").Value = 0.01076420694975565
$ws.Range("P2").Value = 116.2928784986816
$ws.Range("Q2").Value = 176.0177939172235
$ws.Range("A3").Value = "model_30_2_23"
$ws.Range("B3").Value = 0.9998875243994695
$ws.Range("C3").Value = 0.9992155080143623
$ws.Range("D3").Value = 0.9993581415811874
$ws.Range("E3").Value = 0.9998324225448071
$ws.Range("F3").Value = 0.9996487692184403
$ws.Range("G3").Value = 0.000104991035986956
$ws.Range("H3").Value = 0.0007322888333744918
$ws.Range("I3").Value = 0.0004777262360122692
$ws.Range("J3").Value = 0.0001770671487702831
$ws.Range("K3").Value = 0.0003273966923912761
$ws.Range("L3").Value = 0.0006268409962340651
$ws.Range("M3").Value = 0.01024651335757466
$ws.Range("N3").Value = 1.000107976576509
$ws.Range("O3").Value = 0.01068272864441335
$ws.Range("P3").Value = 116.3232711660078
$ws.Range("Q3").Value = 176.0481865845497
$ws.Range("A4").Value = "model_30_2_22"
$ws.Range("B4").Value = 0.9998894235142153
$ws.Range("C4").Value = 0.9992153950460866
$ws.Range("D4").Value = 0.9993687218361061
$ws.Range("E4").Value = 0.9998357634503581
$ws.Range("F4").Value = 0.9996548843618275
$ws.Range("G4").Value = 0.000103218295733317
$ws.Range("H4").Value = 0.0007323942843010761
$ws.Range("I4").Value = 0.0004698515003847899
$ws.Range("J4").Value = 0.0001735370520783023
$ws.Range("K4").Value = 0.0003216965151186955
$ws.Range("L4").Value = 0.000625399550568199
$ws.Range("M4").Value = 0.01015964053169781
$ws.Range("N4").Value = 1.000106153426353
$ws.Range("O4").Value = 0.01059215746248739
$ws.Range("P4").Value = 116.3573288727845
$ws.Range("Q4").Value = 176.0822442913263
$ws.Range("A5").Value = "model_30_2_21"
$ws.Range("B5").Value = 0.9998915275342488
$ws.Range("C5").Value = 0.9992152218884295
$ws.Range("D5").Value = 0.9993805140505242
$ws.Range("E5").Value = 0.9998394427638825
$ws.Range("F5").Value = 0.9996616800982731
$ws.Range("G5").Value = 0.0001012542853877695
$ws.Range("H5").Value = 0.0007325559193732511
$ws.Range("I5").Value = 0.000461074720267812
$ws.Range("J5").Value = 0.0001696493838089981
$ws.Range("K5").Value = 0.0003153619289962143
$ws.Range("L5").Value = 0.0006237373314433825
$ws.Range("M5").Value = 0.01006251883912619
$ws.Range("N5").Value = 1.000104133567121
$ws.Range("O5").Value = 0.01049090109839337
$ws.Range("P5").Value = 116.3957510561116
$ws.Range("Q5").Value = 176.1206664746535
$ws.Range("A6").Value = "model_30_2_20"
$ws.Range("B6").Value = 0.999893754927687
$ws.Range("C6").Value = 0.999215080639991
$ws.Range("D6").Value = 0.9993929634805925
$ws.Range("E6").Value = 0.9998434296603663
$ws.Range("F6").Value = 0.9996689098823838
$ws.Range("G6").Value = 0.00009917511138455049
$ws.Range("H6").Value = 0.0007326877685904993
$ws.Range("I6").Value = 0.0004518087837424171
$ws.Range("J6").Value = 0.000165436714556993
$ws.Range("K6").Value = 0.000308622749149705
$ws.Range("L6").Value = 0.0006219641461297738
$ws.Range("M6").Value = 0.009958670161449796
$ws.Range("N6").Value = 1.000101995269421
$ws.Range("O6").Value = 0.01038263136751189
$ws.Range("P6").Value = 116.4372469369067
$ws.Range("Q6").Value = 176.1621623554485
$ws.Range("A7").Value = "model_30_2_19"
$ws.Range("B7").Value = 0.9998961659519271
$ws.Range("C7").Value = 0.9992148019782789
$ws.Range("D7").Value = 0.9994064822680776
$ws.Range("E7").Value = 0.999847779229896
$ws.Range("F7").Value = 0.9996767722769579
$ws.Range("G7").Value = 0.00009692452608812906
$ws.Range("H7").Value = 0.0007329478870669154
$ws.Range("I7").Value = 0.0004417469394611176
$ws.Range("J7").Value = 0.0001608408345556558
$ws.Range("K7").Value = 0.0003012938870083867
$ws.Range("L7").Value = 0.0006200405972286474
$ws.Range("M7").Value = 0.00984502544883095
$ws.Range("N7").Value = 1.00009968068615
$ws.Range("O7").Value = 0.01026414856419988
$ws.Range("P7").Value = 116.4831559277783
$ws.Range("Q7").Value = 176.2080713463201
$ws.Range("A8").Value = "model_30_2_18"
$ws.Range("B8").Value = 0.999898805699697
$ws.Range("C8").Value = 0.9992145376190139
$ws.Range("D8").Value = 0.999421331025507
$ws.Range("E8").Value = 0.9998525793390866
$ws.Range("F8").Value = 0.9996854209961269
$ws.Range("G8").Value = 0.00009446043741637874
$ws.Range("H8").Value = 0.0007331946548367479
$ws.Range("I8").Value = 0.000430695217167991
$ws.Range("J8").Value = 0.0001557689014176686
$ws.Range("K8").Value = 0.0002932320592928298
$ws.Range("L8").Value = 0.000617803394152442
$ws.Range("M8").Value = 0.009719075954862105
$ws.Range("N8").Value = 1.000097146528291
$ws.Range("O8").Value = 0.01013283713952138
$ws.Range("P8").Value = 116.5346589256044
$ws.Range("Q8").Value = 176.2595743441462
$ws.Range("A9").Value = "model_30_2_17"
$ws.Range("B9").Value = 0.9999015996522825
$ws.Range("C9").Value = 0.9992141986504709
$ws.Range("D9").Value = 0.9994370689701104
$ws.Range("E9").Value = 0.9998577572957578
$ws.Range("F9").Value = 0.9996946369122214
$ws.Range("G9").Value = 0.00009185240531818624
$ws.Range("H9").Value = 0.0007335110670926697
$ws.Range("I9").Value = 0.0004189816853086755
$ws.Range("J9").Value = 0.0001502977237871466
$ws.Range("K9").Value = 0.0002846415239379997
$ws.Range("L9").Value = 0.0006157489379192122
$ws.Range("M9").Value = 0.00958396605368499
$ws.Range("N9").Value = 1.000094464333809
$ws.Range("O9").Value = 0.009991975330135103
$ws.Range("P9").Value = 116.5906551183724
$ws.Range("Q9").Value = 176.3155705369142
$ws.Range("A10").Value = "model_30_2_16"
$ws.Range("B10").Value = 0.9999047185714597
$ws.Range("C10").Value = 0.9992136865986126
$ws.Range("D10").Value = 0.9994548001518077
$ws.Range("E10").Value = 0.9998635006655747
$ws.Range("F10").Value = 0.9997049729631835
$ws.Range("G10").Value = 0.00008894103117109933
$ws.Range("H10").Value = 0.0007339890450257989
$ws.Range("I10").Value = 0.0004057846149828672
$ws.Range("J10").Value = 0.0001442291143990955
$ws.Range("K10").Value = 0.0002750068646909813
$ws.Range("L10").Value = 0.0006134906467024562
$ws.Range("M10").Value = 0.009430855272513694
$ws.Range("N10").Value = 1.000091470171399
$ws.Range("O10").Value = 0.009832346306026335
$ws.Range("P10").Value = 116.6550739578812
$ws.Range("Q10").Value = 176.3799893764231
$ws.Range("A11").Value = "model_30_2_15"
$ws.Range("B11").Value = 0.9999080192277203
$ws.Range("C11").Value = 0.9992130384810661
$ws.Range("D11").Value = 0.9994736548391889
$ws.Range("E11").Value = 0.9998696880453632
$ws.Range("F11").Value = 0.9997160072583665
$ws.Range("G11").Value = 0.00008586001343375642
$ws.Range("H11").Value = 0.0007345940343064164
$ws.Range("I11").Value = 0.0003917513351039708
$ws.Range("J11").Value = 0.0001376913513315647
$ws.Range("K11").Value = 0.0002647213432177678
$ws.Range("L11").Value = 0.000610976115882009
$ws.Range("M11").Value = 0.009266067851778144
$ws.Range("N11").Value = 1.000088301541388
$ws.Range("O11").Value = 0.009660543543633085
$ws.Range("P11").Value = 116.7255846774134
$ws.Range("Q11").Value = 176.4505000959552
$ws.Range("A12").Value = "model_30_2_14"
$ws.Range("B12").Value = 0.9999116453780076
$ws.Range("C12").Value = 0.9992121751214199
$ws.Range("D12").Value = 0.9994945189255499
$ws.Range("E12").Value = 0.9998765148499481
$ws.Range("F12").Value = 0.9997282061850754
$ws.Range("G12").Value = 0.00008247516131021073
$ws.Range("H12").Value = 0.0007353999426391014
$ws.Range("I12").Value = 0.0003762224876932075
$ws.Range("J12").Value = 0.000130477953672172
$ws.Range("K12").Value = 0.0002533502206826897
$ws.Range("L12").Value = 0.0006084172489904764
$ws.Range("M12").Value = 0.009081583634488575
$ws.Range("N12").Value = 1.000084820437113
$ws.Range("O12").Value = 0.009468205451278524
$ws.Range("P12").Value = 116.8060267699671
$ws.Range("Q12").Value = 176.530942188509
$ws.Range("A13").Value = "model_30_2_13"
$ws.Range("B13").Value = 0.9999154184338698
$ws.Range("C13").Value = 0.9992109967409698
$ws.Range("D13").Value = 0.9995164434807058
$ws.Range("E13").Value = 0.9998837183498691
$ws.Range("F13").Value = 0.9997410461686053
$ws.Range("G13").Value = 0.00007895317928093003
$ws.Range("H13").Value = 0.0007364999090643401
$ws.Range("I13").Value = 0.0003599043482033904
$ws.Range("J13").Value = 0.000122866528909116
$ws.Range("K13").Value = 0.0002413815426546377
$ws.Range("L13").Value = 0.0006089062226597745
$ws.Range("M13").Value = 0.00888556015572063
$ws.Range("N13").Value = 1.000081198303485
$ws.Range("O13").Value = 0.009263836847194886
$ws.Range("P13").Value = 116.8933110970479
$ws.Range("Q13").Value = 176.6182265155898
$ws.Range("B14").Value = 0.9999193513228973
$ws.Range("C14").Value = 0.9992096419458248
$ws.Range("D14").Value = 0.9995395040494559
$ws.Range("E14").Value = 0.9998913960529435
$ws.Range("F14").Value = 0.9997546000954891
$ws.Range("G14").Value = 0.00007528200000756083
$ws.Range("H14").Value = 0.0007377645508635038
$ws.Range("I14").Value = 0.0003427406896980615
$ws.Range("J14").Value = 0.0001147540474841344
$ws.Range("K14").Value = 0.000228747368591098
$ws.Range("L14").Value = 0.0006092142697104501
$ws.Range("M14").Value = 0.008676520040175141
$ws.Range("N14").Value = 1.000077422730019
$ws.Range("O14").Value = 0.009045897461157932
$ws.Range("P14").Value = 116.9885389909165
$ws.Range("Q14").Value = 176.7134544094583
$ws.Range("A15").Value = "model_30_2_11"
$ws.Range("B15").Value = 0.9999235172956782
$ws.Range("C15").Value = 0.9992079744632271
$ws.Range("D15").Value = 0.9995642209943795
$ws.Range("E15").Value = 0.9998996276722
$ws.Range("F15").Value = 0.9997691334303647
$ws.Range("G15").Value = 0.00007139324728172152
$ws.Range("H15").Value = 0.0007393210726745164
$ws.Range("I15").Value = 0.0003243442135936745
$ws.Range("J15").Value = 0.0001060562823234697
$ws.Range("K15").Value = 0.0002152002479585721
$ws.Range("L15").Value = 0.0006090273374019322
$ws.Range("M15").Value = 0.008449452484138928
$ws.Range("N15").Value = 1.000073423396149
$ws.Range("O15").Value = 0.008809163169166618
$ws.Range("P15").Value = 117.0946145379551
$ws.Range("Q15").Value = 176.819529956497
$ws.Range("A16").Value = "model_30_2_10"
$ws.Range("B16").Value = 0.9999279508255102
$ws.Range("C16").Value = 0.9992059175095817
$ws.Range("D16").Value = 0.9995908313445283
$ws.Range("E16").Value = 0.9999084932935468
$ws.Range("F16").Value = 0.9997847780675528
$ws.Range("G16").Value = 0.00006725474179299524
$ws.Range("H16").Value = 0.0007412411486127047
$ws.Range("I16").Value = 0.0003045385024851572
$ws.Range("J16").Value = 0.00009668861235765237
$ws.Range("K16").Value = 0.0002006172366224272
$ws.Range("L16").Value = 0.0006076395132619014
$ws.Range("M16").Value = 0.008200898840553713
$ws.Range("N16").Value = 1.00006916720751
$ws.Range("O16").Value = 0.008550028082397007
$ws.Range("P16").Value = 117.21404606407
$ws.Range("Q16").Value = 176.9389614826119
$ws.Range("A17").Value = "model_30_2_9"
$ws.Range("B17").Value = 0.9999323744989834
$ws.Range("C17").Value = 0.9992030899608439
$ws.Range("D17").Value = 0.9996180753419184
$ws.Range("E17").Value = 0.9999174887139486
$ws.Range("F17").Value = 0.9998007571395686
$ws.Range("G17").Value = 0.00006312543678255215
$ws.Range("H17").Value = 0.0007438805412442581
$ws.Range("I17").Value = 0.000284261176605135
$ws.Range("J17").Value = 0.00008718379298499772
$ws.Range("K17").Value = 0.0001857224847950663
$ws.Range("L17").Value = 0.000606734565526585
$ws.Range("M17").Value = 0.007945151778446536
$ws.Range("N17").Value = 1.000064920480976
$ws.Range("O17").Value = 0.008283393338386524
$ws.Range("P17").Value = 117.3407735021303
$ws.Range("Q17").Value = 177.0656889206721
$ws.Range("A18").Value = "model_30_2_8"
$ws.Range("B18").Value = 0.9999367958654202
$ws.Range("C18").Value = 0.9991999989111137
$ws.Range("D18").Value = 0.9996456798041765
$ws.Range("E18").Value = 0.999926813733127
$ws.Range("F18").Value = 0.9998170665774955
$ws.Range("G18").Value = 0.00005899828528936955
$ws.Range("H18").Value = 0.0007467659004861578
$ws.Range("I18").Value = 0.0002637155617698407
$ws.Range("J18").Value = 0.0000773307100851991
$ws.Range("K18").Value = 0.0001705197852814585
$ws.Range("L18").Value = 0.0006069965949203784
$ws.Range("M18").Value = 0.007681034128902797
$ws.Range("N18").Value = 1.000060675969197
$ws.Range("O18").Value = 0.008008031653702802
$ws.Range("P18").Value = 117.4760043547459
$ws.Range("Q18").Value = 177.2009197732877
$ws.Range("A19").Value = "model_30_2_7"
$ws.Range("B19").Value = 0.9999413468038079
$ws.Range("C19").Value = 0.9991959636344911
$ws.Range("D19").Value = 0.9996749683740933
$ws.Range("E19").Value = 0.9999363570655416
$ws.Range("F19").Value = 0.999834164914763
$ws.Range("G19").Value = 0.00005475018406761682
$ws.Range("H19").Value = 0.0007505326540850263
$ws.Range("I19").Value = 0.0002419164891793069
$ws.Range("J19").Value = 0.00006724695115428183
$ws.Range("K19").Value = 0.0001545817201667944
$ws.Range("L19").Value = 0.0006058170328005303
$ws.Range("M19").Value = 0.0073993367316008
$ws.Range("N19").Value = 1.000056307068344
$ws.Range("O19").Value = 0.007714341815003658
$ws.Range("P19").Value = 117.6254596546154
$ws.Range("Q19").Value = 177.3503750731572
$ws.Range("A20").Value = "model_30_2_6"
$ws.Range("B20").Value = 0.9999460450343235
$ws.Range("C20").Value = 0.9991909126329679
$ws.Range("D20").Value = 0.9997060832112992
$ws.Range("E20").Value = 0.9999460413905044
$ws.Range("F20").Value = 0.9998520758666837
$ws.Range("G20").Value = 0.00005036459210973895
$ws.Range("H20").Value = 0.0007552475422937247
$ws.Range("I20").Value = 0.0002187581514106855
$ws.Range("J20").Value = 0.00005701421545038958
$ws.Range("K20").Value = 0.0001378861834305375
$ws.Range("L20").Value = 0.0006014410272122276
$ws.Range("M20").Value = 0.007096801540816746
$ws.Range("N20").Value = 1.000051796767049
$ws.Range("O20").Value = 0.007398927074813746
$ws.Range("P20").Value = 117.7924443345745
$ws.Range("Q20").Value = 177.5173597531163
$ws.Range("A21").Value = "model_30_2_5"
$ws.Range("B21").Value = 0.9999506281450994
$ws.Range("C21").Value = 0.9991863819508725
$ws.Range("D21").Value = 0.99973495209348
$ws.Range("E21").Value = 0.9999559109811018
$ws.Range("F21").Value = 0.9998691951596518
$ws.Range("G21").Value = 0.0000460864593758053
$ws.Range("H21").Value = 0.0007594767351557761
$ws.Range("I21").Value = 0.0001972714465271153
$ws.Range("J21").Value = 0.0000465857227596365
$ws.Range("K21").Value = 0.0001219285846433759
$ws.Range("L21").Value = 0.0005880847759007707
$ws.Range("M21").Value = 0.006788700860680584
$ws.Range("N21").Value = 1.000047396980705
$ws.Range("O21").Value = 0.007077709910867856
$ws.Range("P21").Value = 117.9699827479485
$ws.Range("Q21").Value = 177.6948981664903
$ws.Range("A22").Value = "model_30_2_1"
$ws.Range("B22").Value = 0.9999640593107416
$ws.Range("C22").Value = 0.9991819568023446
$ws.Range("D22").Value = 0.9998348883457653
$ws.Range("E22").Value = 0.9999825404646437
$ws.Range("F22").Value = 0.999924186091829
$ws.Range("G22").Value = 0.00003354905580878713
$ws.Range("H22").Value = 0.0007636074170650123
$ws.Range("I22").Value = 0.0001228902929172943
$ws.Range("J22").Value = 0.00001844824616073892
$ws.Range("K22").Value = 0.00007066926953901661
$ws.Range("L22").Value = 0.0005577189059106345
$ws.Range("M22").Value = 0.0057921546775606
$ws.Range("N22").Value = 1.000034503061688
$ws.Range("O22").Value = 0.006038738693597466
$ws.Range("P22").Value = 118.6050036755015
$ws.Range("Q22").Value = 178.3299190940433
$ws.Range("A23").Value = "model_30_2_4"
$ws.Range("B23").Value = 0.9999549916888935
$ws.Range("C23").Value = 0.9991803040811595
$ws.Range("D23").Value = 0.9997648993297821
$ws.Range("E23").Value = 0.9999642790544758
$ws.Range("F23").Value = 0.9998858939464286
$ws.Range("G23").Value = 0.00004201328278145985
$ws.Range("H23").Value = 0.0007651501597451215
$ws.Range("I23").Value = 0.0001749821377664005
$ws.Range("J23").Value = 0.00003774377626196198
$ws.Range("K23").Value = 0.0001063629570141812
$ws.Range("L23").Value = 0.000577905230009324
$ws.Range("M23").Value = 0.00648176540623462
$ws.Range("N23").Value = 1.000043207978662
$ws.Range("O23").Value = 0.006757707578682145
$ws.Range("P23").Value = 118.1550494659562
$ws.Range("Q23").Value = 177.879964884498
$ws.Range("A24").Value = "model_30_2_0"
$ws.Range("B24").Value = 0.9999647141514889
$ws.Range("C24").Value = 0.9991801440888387
$ws.Range("D24").Value = 0.9998566185784571
$ws.Range("E24").Value = 0.9999856306298913
$ws.Range("F24").Value = 0.9999346129533665
$ws.Range("G24").Value = 0.00003293779071537396
$ws.Range("H24").Value = 0.0007652995055537825
$ws.Range("I24").Value = 0.0001067167849172903
$ws.Range("J24").Value = 0.00001518308886980715
$ws.Range("K24").Value = 0.00006094996201062589
$ws.Range("L24").Value = 0.0005704697083375326
$ws.Range("M24").Value = 0.005739145469089798
$ws.Range("N24").Value = 1.000033874414571
$ws.Range("O24").Value = 0.005983472773377874
$ws.Range("P24").Value = 118.6417798108569
$ws.Range("Q24").Value = 178.3666952293987
$ws.Range("A25").Value = "model_30_2_2"
$ws.Range("B25").Value = 0.9999628402321342
$ws.Range("C25").Value = 0.9991786534397771
$ws.Range("D25").Value = 0.9998175489617455
$ws.Range("E25").Value = 0.9999781363162995
$ws.Range("F25").Value = 0.9999147674491407
$ws.Range("G25").Value = 0.00003468701217737945
$ws.Range("H25").Value = 0.0007666909610209288
$ws.Range("I25").Value = 0.0001357957537163827
$ws.Range("J25").Value = 0.00002310179570395633
$ws.Range("K25").Value = 0.0000794487747101695
$ws.Range("L25").Value = 0.0005524603410412244
$ws.Range("M25").Value = 0.00588956808071521
$ws.Range("N25").Value = 1.000035673377151
$ws.Range("O25").Value = 0.006140299186997913
$ws.Range("P25").Value = 118.5382904598227
$ws.Range("Q25").Value = 178.2632058783645
$ws.Range("A26").Value = "model_30_2_3"
$ws.Range("B26").Value = 0.999958944816002
$ws.Range("C26").Value = 0.9991735606623796
$ws.Range("D26").Value = 0.9997945704667638
$ws.Range("E26").Value = 0.9999709600483068
$ws.Range("F26").Value = 0.9999015263102948
$ws.Range("G26").Value = 0.00003832321214782787
$ws.Range("H26").Value = 0.0007714448451743175
$ws.Range("I26").Value = 0.0001528983258648395
$ws.Range("J26").Value = 0.00003068444643006106
$ws.Range("K26").Value = 0.00009179138614745031
$ws.Range("L26").Value = 0.0005701054632381329
$ws.Range("M26").Value = 0.006190574460244208
$ws.Range("N26").Value = 1.000039412976638
$ws.Range("O26").Value = 0.006454120031272564
$ws.Range("P26").Value = 118.3389095681153
$ws.Range("Q26").Value = 178.0638249866572
